# Update data map for Italy: add ITA1 region/ageband rows to Sheet1
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("A15").Value = "ITA1"
$ws.Range("B15").Value = "region"
$ws.Range("C15").Value = "data/derived/ITA/ITA_regions.RDS"
$ws.Range("D15").Value = "marginal"
$ws.Range("E15").Value = "aggregate"

$ws.Range("A16").Value = "ITA1"
$ws.Range("B16").Value = "ageband"
$ws.Range("C16").Value = "data/derived/ITA/ITA_agebands.RDS"
$ws.Range("D16").Value = "marginal"
$ws.Range("E16").Value = "aggregate"
